$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.005438786316491882
$ws.Range("B3").Value = 0.05968737632435654
$ws.Range("B4").Value = 0.02419491912793581
$ws.Range("B5").Value = 0.03621861882594592
$ws.Range("B6").Value = 0.1554086597736452
$ws.Range("B7").Value = 0.05801805089369379
$ws.Range("B8").Value = 0.1230389365061619
$ws.Range("B9").Value = 0.03967146056162398
$ws.Range("B10").Value = 0.1839243143310103
$ws.Range("B11").Value = 0.01897628365072135
$ws.Range("B12").Value = 0.03493978196040689
$ws.Range("B13").Value = 0.02546352325793216
$ws.Range("B14").Value = 0.02533537266980946
$ws.Range("B15").Value = 0.02346693498438158
$ws.Range("B16").Value = 0.03139293223836195
$ws.Range("B17").Value = 0.02559912392072757
$ws.Range("B18").Value = 0.006638500410475031
$ws.Range("B19").Value = 0.02814984173776284
$ws.Range("B20").Value = 0.02522989868018976
$ws.Range("B21").Value = 0.005011118523811611
$ws.Range("B22").Value = 0.01568976138265739
$ws.Range("B23").Value = 0.006219801865620575
$ws.Range("B24").Value = 0.0154825700666057
$ws.Range("B25").Value = 0.0102661265955363
$ws.Range("B26").Value = 0.01653730539413442
